$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper cell (far outside the used range) used to build clean text values
# for the date column without Excel auto-converting the "yyyy-mm-dd"-looking
# text into a real date (and thereby attaching a new/extra number format
# style to the cell). TRIM() forces the formula result to be plain text,
# and pasting its *value* onto the destination cell keeps it as a literal
# shared string with the worksheet's default (unstyled) cell format --
# matching how the existing date cells (A2:A68) are stored.
$helper = $ws.Range("Z1")

# Row 69: raw/clean SSA data for 2020-08-07
$helper.Formula = "=TRIM(""2020-08-07"")"
$helper.Copy() | Out-Null
$ws.Range("A69").PasteSpecial(-4163) | Out-Null

$ws.Range("B69").Value = 469407
$ws.Range("C69").Value = 513144
$ws.Range("D69").Value = 89155
$ws.Range("E69").Value = 51311
$ws.Range("F69").Value = 26.74

# Row 70: raw/clean SSA data for 2020-08-08
$helper.Formula = "=TRIM(""2020-08-08"")"
$helper.Copy() | Out-Null
$ws.Range("A70").PasteSpecial(-4163) | Out-Null

$ws.Range("B70").Value = 475902
$ws.Range("C70").Value = 520970
$ws.Range("D70").Value = 89025
$ws.Range("E70").Value = 52006
$ws.Range("F70").Value = 26.67

$helper.Clear() | Out-Null
$excel.CutCopyMode = $false
